# Updated to use late january bvc sheet
# Fixed first bvc store not being loaded correctly
#
# Appends one new beer-code lookup row (short code -> full name -> category)
# to the bottom of the (single) sheet, just after the current last row (85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 86: Appl Crisp -> Apple Crisp Porter -> Other
$newRow = 86
$ws.Cells.Item($newRow, 1).Value = "Appl Crisp"
$ws.Cells.Item($newRow, 2).Value = "Apple Crisp Porter"
$ws.Cells.Item($newRow, 3).Value = "Other"

# Bring the freshly added row into view / selection, matching the sheet
# having been scrolled down and the new (blank) row below it selected
# while editing.
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows($newRow + 3).Select()
